$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data rows (rows 2:89), keeping the header row (row 1)
# intact with its existing style/formatting.
$ws.Rows("2:89").Clear()

# New "Active cases" snapshot data: (Cluster name, Active cases, Exist)
$data = @(
    @("21 Moore Street Apartment Complex Moonee Ponds", 5, "old"),
    @("A Treat of France Cafe Carlton North", 5, "old"),
    @("A1 Bakery Brunswick", 6, "new"),
    @("Acquire BPO Southbank", 36, "new"),
    @("Acquire BPO Southbank", 38, "old"),
    @("Al-Taqwa College Truganina", 28, "new"),
    @("Al-Taqwa College Truganina", 30, "old"),
    @("Amiga Montessori Craigieburn", 19, "old"),
    @("Amiga Montessori Craigieburn", 22, "new"),
    @("Broadmeadows Medical Centre Broadmeadows", 5, "new"),
    @("CS Square Caroline Springs", 5, "new"),
    @("CS Square Caroline Springs", 8, "old"),
    @("Campbellfield Heights Primary School", 6, "new"),
    @("Campbellfield Heights Primary School Campbellfield", 6, "old"),
    @("Cedars Medical Clinic Coburg", 39, "old"),
    @("Cedars Medical Clinic Coburg", 44, "new"),
    @("Chemist Warehouse Campbellfield DC", 5, "old"),
    @("City of Greater Shepparton Community", 6, "new"),
    @("City of Greater Shepparton Community", 8, "old"),
    @("City of Hobsons Bay Community", 10, "old"),
    @("City of Hobsons Bay Community", 12, "new"),
    @("City of Moreland Social Network", 5, "old"),
    @("City of Wyndham Community", 7, "new"),
    @("City of Wyndham Community", 8, "old"),
    @("Classy Cabinets and Kitchens Craigieburn", 16, "new"),
    @("Classy Cabinets and Kitchens Craigieburn", 18, "old"),
    @("Coles Broadmeadows Central Shopping Centre", 9, "old"),
    @("Coles Broadmeadows Central Shopping Centre", 11, "new"),
    @("Coles Campbellfield Plaza Campbellfield", 12, "new"),
    @("Coles Coburg North Village Aug", 7, "old"),
    @("Coles Coburg North Village August", 8, "new"),
    @("Coles Greenvale Shopping Centre", 5, "old"),
    @("Coles Greenvale Shopping Centre", 6, "new"),
    @("Coles Pakenham Place Shopping Centre", 6, "new"),
    @("Crossroads Logistics Sunshine North", 6, "new"),
    @("Crossroads Logistics Sunshine North", 7, "old"),
    @("Elite Smart Community Care Campbellfield", 5, "new"),
    @("Ernst and Young 8 Exhibition Street Melbourne", 12, "old"),
    @("Fitzroy Community School Fitzroy North", 48, "old"),
    @("Fitzroy Community School Fitzroy North", 50, "new"),
    @("IGA Meadow Heights Shopping Centre Meadow Heights", 5, "old"),
    @("IGA Meadow Heights Shopping Centre Meadow Heights", 7, "new"),
    @("Islamic College of Melbourne Tarneit", 6, "old"),
    @("Islamic College of Melbourne Tarneit", 7, "new"),
    @("Kasr Sweets Coolaroo", 22, "new"),
    @("Kasr Sweets Coolaroo", 25, "old"),
    @("Kids House Early Learning Cheltenham", 6, "old"),
    @("Kids House Early Learning Cheltenham", 10, "new"),
    @("Learning Nest Early Learning Centre Meadow Heights", 9, "old"),
    @("Learning Nest Early Learning Centre Meadow Heights", 10, "new"),
    @("Louis Vuitton Melbourne Crown Southbank", 5, "new"),
    @("Melbourne Metropolitan Remand Centre Ravenhall", 5, "new"),
    @("Melbourne Truck Repairs Campbellfield", 5, "new"),
    @("Mercy Hospital for Women Heidelberg", 6, "new"),
    @("Metricon Homes Lot 1815 Canning Dr Mickleham", 5, "old"),
    @("MyCentre Childcare Broadmeadows", 19, "new"),
    @("MyCentre Childcare Broadmeadows", 25, "old"),
    @("Northern Health Northern Hospital Epping Emergency Department Tier 1B", 39, "old"),
    @("Northern Health Northern Hospital Epping Emergency Department Tier 1B", 42, "new"),
    @("OnQ Plumbing and Excavations Cragieburn", 7, "old"),
    @("OnQ Plumbing and Excavations Craigieburn", 8, "new"),
    @("Oporto Coolaroo", 13, "old"),
    @("Oporto Coolaroo", 14, "new"),
    @("Paisley Park Early Learning Centre Bundoora", 7, "old"),
    @("Paisley Park Early Learning Centre Bundoora", 8, "new"),
    @("Panorama Construction Site Whitehorse Rd Box Hill", 67, "new"),
    @("Panorama Construction Site Whitehorse Rd Box Hill", 76, "old"),
    @("Pascoe Vale Primary School Pascoe Vale", 6, "new"),
    @("Pascoe Vale Primary School Pascoe Vale", 7, "old"),
    @("St Vincents Hospital Emergency Department Melbourne", 6, "old"),
    @("St Vincents Hospital Emergency Department Melbourne", 7, "new"),
    @("Tek Foods Somerton", 7, "old"),
    @("Tek Foods Somerton", 9, "new"),
    @("Tip Top Warehouse Dandenong", 9, "old"),
    @("Tip Top Warehouse Dandenong", 10, "new"),
    @("Tunstall Fresh Tunstall Square Shopping Centre Doncaster East", 5, "new"),
    @("Turnstall Fresh Turnstall Square Shopping Centre Donaster East", 5, "old"),
    @("Western Health Footscray Hospital Emergency Department", 6, "new"),
    @("Western Health Footscray Hospital Emergency Department", 7, "old"),
    @("Yara Childcare Centre Truganina", 30, "new"),
    @("Yara Childcare Centre Truganina", 31, "old")

)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = $row[2]
    $r++
}
